$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '62.428.74'
Set-TextValue 'E2' '  +2.59%  '
Set-TextValue 'D3' '2.427.77'
Set-TextValue 'E3' '  +3.43%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '556.56'
Set-TextValue 'E5' '  +2.29%  '
Set-TextValue 'D6' '143.97'
Set-TextValue 'E6' '  +5.73%  '
Set-TextValue 'E8' '  +2.17%  '
Set-TextValue 'D9' '2.429.98'
Set-TextValue 'E9' '  +3.57%  '
Set-TextValue 'E10' '  +5.62%  '
Set-TextValue 'E11' '  +0.37%  '
Set-TextValue 'E12' '  +1.99%  '
Set-TextValue 'E13' '  +4.36%  '
Set-TextValue 'D14' '26.37'
Set-TextValue 'E14' '  +7.52%  '
Set-TextValue 'E15' '  +10.37%  '
Set-TextValue 'D16' '2.866.30'
Set-TextValue 'E16' '  +3.47%  '
Set-TextValue 'D17' '62.257.15'
Set-TextValue 'E17' '  +2.63%  '
Set-TextValue 'D18' '2.427.97'
Set-TextValue 'E18' '  +3.60%  '
Set-TextValue 'D19' '11.12'
Set-TextValue 'E19' '  +5.03%  '
Set-TextValue 'D20' '325.36'
Set-TextValue 'E20' '  +2.07%  '
Set-TextValue 'E21' '  +2.00%  '
Set-TextValue 'E22' '  +4.10%  '
Set-TextValue 'E23' '  +0.12%  '
Set-TextValue 'D24' '1.78'
Set-TextValue 'E24' '  +3.64%  '
Set-TextValue 'D25' '65.20'
Set-TextValue 'E25' '  +3.41%  '
Set-TextValue 'D26' '9.16'
Set-TextValue 'E26' '  +11.60%  '
Set-TextValue 'D27' '573.27'
Set-TextValue 'E27' '  +15.71%  '
Set-TextValue 'D28' '2.538.35'
Set-TextValue 'E28' '  +3.23%  '
Set-TextValue 'E29' '  +0.32%  '
Set-TextValue 'D30' '0.0₃0946'
Set-TextValue 'E30' '  +10.49%  '
Set-TextValue 'D31' '8.41'
Set-TextValue 'E31' '  +6.50%  '
Set-TextValue 'D32' '1.45'
Set-TextValue 'E32' '  +6.44%  '
Set-TextValue 'E33' '  +2.03%  '
Set-TextValue 'D34' '1.87'
Set-TextValue 'E34' '  +4.89%  '
Set-TextValue 'E35' '  +5.15%  '
Set-TextValue 'E36' '  +9.67%  '
Set-TextValue 'B37' 'FirstDigitalUSD'
Set-TextValue 'C37' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D37' '0.999'
Set-TextValue 'E37' '  +0.08%  '
Set-TextValue 'B38' 'NEARProtocol'
Set-TextValue 'C38' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D38' '4.83'
Set-TextValue 'E38' '  +5.81%  '
Set-TextValue 'B39' 'PolygonEcosystemToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D39' '0.386'
Set-TextValue 'E39' '  +2.79%  '
Set-TextValue 'B40' 'Stacks'
Set-TextValue 'C40' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D40' '1.90'
Set-TextValue 'E40' '  +5.16%  '
Set-TextValue 'D41' '18.89'
Set-TextValue 'E41' '  +1.78%  '
Set-TextValue 'D42' '147.50'
Set-TextValue 'E42' '  +4.26%  '
Set-TextValue 'D44' '41.67'
Set-TextValue 'E44' '  +2.73%  '
Set-TextValue 'D45' '2.32'
Set-TextValue 'E45' '  +14.17%  '
Set-TextValue 'D46' '152.20'
Set-TextValue 'E46' '  +7.31%  '
Set-TextValue 'E47' '  +3.05%  '
Set-TextValue 'E48' '  +6.16%  '
Set-TextValue 'E49' '  +8.62%  '
Set-TextValue 'E50' '  +9.70%  '
Set-TextValue 'E51' '  +4.62%  '
